$wb = $excel.ActiveWorkbook

# Rename the sheet tabs from spaced names to underscored survey-code-style names.
$wb.Worksheets.Item(1).Name = "Test_Yearly"
$wb.Worksheets.Item(2).Name = "Test_Weekly"

# The bold/"applyFont" cell style (cellXfs index 2) is no longer used on the
# header-row labels in either sheet; those cells revert to the default style.
$ws1 = $wb.Worksheets.Item("Test_Yearly")
$ws1.Range("A5").Style = "Normal"
$ws1.Range("A6").Style = "Normal"
$ws1.Range("A7").Style = "Normal"
$ws1.Range("K4").Style = "Normal"

$ws2 = $wb.Worksheets.Item("Test_Weekly")
$ws2.Range("A5").Style = "Normal"
$ws2.Range("A6").Style = "Normal"
$ws2.Range("A7").Style = "Normal"
$ws2.Range("K4").Style = "Normal"
